$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (calibrated) data for rows 2-12, columns A-D
# Numeric literals are kept as strings and cast to [double] because the
# script engine's parser does not support scientific-notation number
# literals directly.
$data = @(
    @("56898.967185", "-7.402089929e-06",       "-1.1925646062e-05",      "-2.5329835517e-05"),
    @("56916.699185", "-4.7993912143e-05",      "-0.00010193708406",      "-9.9468749202e-05"),
    @("56927.699186", "-5.6037888207e-05",      "-0.00023315827516",      "-0.00014170938367"),
    @("56938.567186", "-7.5421499197e-05",      "-0.00036139814919",      "-0.00019064895824"),
    @("56949.435187", "-0.0001339691",          "-0.0004897914",          "-0.0002453469"),
    @("56960.367188", "-0.0002014021",          "-0.0006097965",          "-0.0003018801"),
    @("56978.167188", "-0.0001424692",          "-0.0004881551",          "-0.0002451157"),
    @("56990.831189", "-7.901043384099999e-05", "-0.00035218664206",      "-0.00020052239348"),
    @("57001.231189", "-5.5720426976e-05",      "-0.00022620906777",      "-0.00014607513843"),
    @("57011.56719",  "-2.4025626554e-05",      "-7.608194195399999e-05", "-0.00010670966426"),
    @("57021.63119",  "-7.2268442087e-06",      "-1.1698795606e-05",      "-2.495058565e-05")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = [double]$values[0]
    $ws.Cells.Item($row, 2).Value = [double]$values[1]
    $ws.Cells.Item($row, 3).Value = [double]$values[2]
    $ws.Cells.Item($row, 4).Value = [double]$values[3]
}
